$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 06:52"

# --- Simple in-place numeric updates --------------------------------------
# Row 28 (Singapur)
$ws.Range("F28").Value = 21

# Row 63 (Tailandia)
$ws.Range("B63").Value = 2960
$ws.Range("C63").Value = 6
$ws.Range("D63").Value = 2719
$ws.Range("E63").Value = 187

# Row 97 (Kirguistan)
$ws.Range("B97").Value = 756
$ws.Range("C97").Value = 10
$ws.Range("D97").Value = 504
$ws.Range("E97").Value = 244

# --- El Salvador overtakes Jamaica / Reunion / Kenia in the ranking ------
# (sorted descending by "Casos totales"); rows 119-122 shuffle down and
# El Salvador's updated totals land in row 119.
$ws.Range("A119").Value = "El Salvador"
$ws.Range("B119").Value = 424
$ws.Range("C119").Value = 29
$ws.Range("D119").Value = 124
$ws.Range("E119").Value = 290
$ws.Range("F119").Value = 2
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 10

$ws.Range("A120").Value = "Jamaica"
$ws.Range("B120").Value = 422
$ws.Range("C120").Value = 26
$ws.Range("D120").Value = 29
$ws.Range("E120").Value = 385
$ws.Range("F120").Value = 3
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 8

$ws.Range("A121").Value = "Reunion"
$ws.Range("B121").Value = 420
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 300
$ws.Range("E121").Value = 120
$ws.Range("F121").Value = 2
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 0

$ws.Range("A122").Value = "Kenia"
$ws.Range("B122").Value = 396
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 144
$ws.Range("E122").Value = 235
$ws.Range("F122").Value = 2
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 17

# A123 (Estado de Palestina) and its row stay as-is.
